# Update column G ("K") values in the active worksheet.
# The source data previously stored a raw strike count ("Strike#");
# it has been regenerated to store K (computed from std/mean of the
# underlying s_vals) instead. Only column G changes; every other
# column is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2"  = 0
    "G3"  = 0
    "G4"  = 1
    "G5"  = 2
    "G6"  = 0
    "G7"  = 1
    "G8"  = 1
    "G9"  = 1
    "G10" = 1
    "G11" = 0
    "G12" = 1
    "G13" = 0
    "G14" = 0
    "G15" = 1
    "G17" = 1
    "G18" = 1
    "G19" = 0
    "G20" = 0
    "G21" = 3
    "G22" = 3
    "G23" = 1
    "G24" = 1
    "G25" = 2
    "G26" = 0
    "G27" = 0
    "G28" = 0
    "G29" = 2
    "G30" = 0
    "G31" = 3
    "G32" = 4
    "G33" = 0
    "G34" = 2
    "G35" = 1
    "G36" = 0
    "G37" = 0
    "G38" = 1
    "G39" = 2
    "G40" = 1
    "G41" = 1
    "G42" = 1
    "G43" = 0
    "G44" = 1
    "G45" = 1
    "G46" = 0
    "G47" = 1
    "G48" = 0
    "G49" = 2
    "G50" = 0
    "G51" = 0
    "G52" = 0
    "G53" = 2
    "G54" = 1
    "G55" = 3
    "G56" = 1
    "G57" = 1
    "G58" = 1
    "G59" = 0
    "G60" = 0
    "G61" = 0
    "G62" = 1
    "G63" = 0
    "G64" = 3
    "G65" = 4
    "G66" = 4
    "G67" = 1
    "G68" = 2
    "G69" = 0
    "G70" = 1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
